# Update scripts with new TPM-derived values.
# The "ECs" sending-cluster rows are dropped, and the "MuSCs" sending-cluster
# rows (formerly rows 5-7) move up to become rows 2-4 with refreshed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the three "ECs" sender rows (old rows 2,3,4); this shifts the old
# "MuSCs" sender rows (5,6,7) up to become rows 2,3,4.
$ws.Range("A2:A4").EntireRow.Delete()

# Row 2: MuSCs -> Cdh1 -> Ptprm -> ECs
$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(2, 13).Value = 77.07050600000001
$ws.Cells.Item(2, 14).Value = 231.211518
$ws.Cells.Item(2, 15).Value = 0.7967384544746254
$ws.Cells.Item(2, 16).Value = 0.7967384544746255
$ws.Cells.Item(2, 17).Value = 73.04100304463334
$ws.Cells.Item(2, 18).Value = 657.3690274017
$ws.Cells.Item(2, 19).Value = 0.7967384544746254
$ws.Cells.Item(2, 20).Value = 0.7967384544746255

# Row 3: MuSCs -> Cdh1 -> Ptprm -> FAPs
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(3, 15).Value = 0.1548701728536178
$ws.Cells.Item(3, 16).Value = 0.1548701728536178
$ws.Cells.Item(3, 17).Value = 14.19772411309444
$ws.Cells.Item(3, 18).Value = 127.77951701785
$ws.Cells.Item(3, 19).Value = 0.1548701728536178
$ws.Cells.Item(3, 20).Value = 0.1548701728536178

# Row 4: MuSCs -> Cdh1 -> Ptprm -> MuSCs
$ws.Cells.Item(4, 9).Value = 1
$ws.Cells.Item(4, 10).Value = 1
$ws.Cells.Item(4, 15).Value = 0.04839137267175684
$ws.Cells.Item(4, 16).Value = 0.04839137267175685
$ws.Cells.Item(4, 17).Value = 4.436279407377778
$ws.Cells.Item(4, 18).Value = 39.9265146664
$ws.Cells.Item(4, 19).Value = 0.04839137267175684
$ws.Cells.Item(4, 20).Value = 0.04839137267175685
